$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numbers formatted with "." as both thousands
# and decimal separators (e.g. "67.553.24"), stored as plain text in the
# workbook. Force text format while writing so COM does not coerce
# decimal-looking values (e.g. "6.17") into real numbers, then restore
# the original style so no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$origStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '67.553.24'
$ws.Range("E2").Value = '  -2.81%  '
$ws.Range("D3").Value = '3.725.25'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '592.43'
$ws.Range("E5").Value = '  -3.26%  '
$ws.Range("E6").Value = '  -3.39%  '
$ws.Range("D7").Value = '3.727.24'
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("E10").Value = '  -5.24%  '
$ws.Range("D11").Value = '6.17'
$ws.Range("E11").Value = '  -5.96%  '
$ws.Range("E12").Value = '  -4.51%  '
$ws.Range("D13").Value = '37.42'
$ws.Range("E13").Value = '  -5.86%  '
$ws.Range("E14").Value = '  -4.66%  '
$ws.Range("D15").Value = '4.350.11'
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '3.725.85'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").Value = '67.542.98'
$ws.Range("E17").Value = '  -2.96%  '
$ws.Range("E18").Value = '  -5.05%  '
$ws.Range("D19").Value = '7.12'
$ws.Range("E19").Value = '  -3.89%  '
$ws.Range("D20").Value = '16.12'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").Value = '487.27'
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").Value = '8.97'
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("D23").Value = '0.713'
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").Value = '83.33'
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("E25").Value = '  -9.21%  '
$ws.Range("D26").Value = '0.0000141'
$ws.Range("E26").Value = '  +6.41%  '
$ws.Range("D27").Value = '12.10'
$ws.Range("E27").Value = '  -5.45%  '
$ws.Range("D28").Value = '10.19'
$ws.Range("E28").Value = '  -7.02%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").Value = '2.92'
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D33").Value = '7.62'
$ws.Range("E33").Value = '  -4.54%  '
$ws.Range("E34").Value = '  -4.40%  '
$ws.Range("D36").Value = '0.995'
$ws.Range("E36").Value = '  -4.90%  '
$ws.Range("E37").Value = '  -2.17%  '
$ws.Range("D38").Value = '5.68'
$ws.Range("E38").Value = '  -6.46%  '
$ws.Range("D39").Value = '0.323'
$ws.Range("E39").Value = '  -6.35%  '
$ws.Range("D40").Value = '447.03'
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("D41").Value = '48.79'
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("E42").Value = '  -3.50%  '
$ws.Range("D43").Value = '2.84'
$ws.Range("E43").Value = '  -6.39%  '
$ws.Range("D44").Value = '8.24'
$ws.Range("E44").Value = '  -3.09%  '
$ws.Range("D45").Value = '41.22'
$ws.Range("E45").Value = '  -6.82%  '
$ws.Range("D46").Value = '140.83'
$ws.Range("E46").Value = '  +1.77%  '
$ws.Range("D47").Value = '2.782.44'
$ws.Range("E47").Value = '  -5.46%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").Value = '0.0346'
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").Value = '25.56'
$ws.Range("E50").Value = '  -4.85%  '
$ws.Range("D51").Value = '23.17'
$ws.Range("E51").Value = '  +8.79%  '
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "2.37"
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "32.18"
$ws.Range("E32").Value = "  +6.25%  "

$priceRange.Style = $origStyle
